$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (rich-text runs inside shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "12"
$ws.Range("C9").Characters(27, 9).Text = "3/17/2025"
$ws.Range("C9").Characters(47, 9).Text = "3/23/2025"

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -50
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 133.333333333333
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 19
$ws.Range("K16").Value = 42.105263157894
$ws.Range("L16").Value = -12.903225806451
$ws.Range("M16").Value = 3.846153846153
$ws.Range("N16").Value = -84.023668639053
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -38.095238095238
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 54
$ws.Range("K17").Value = -20.370370370370
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 48.275862068965
$ws.Range("N17").Value = 30.303030303030
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -65
$ws.Range("I18").Value = 28
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = 7.692307692307
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = -61.111111111111
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -13.793103448275
$ws.Range("I19").Value = 79
$ws.Range("J19").Value = 81
$ws.Range("K19").Value = -2.469135802469
$ws.Range("L19").Value = -30.088495575221
$ws.Range("M19").Value = 64.583333333333
$ws.Range("N19").Value = -13.186813186813
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -66.666666666666
$ws.Range("N20").Value = -97.468354430379
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -36
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -24.050632911392
$ws.Range("I21").Value = 182
$ws.Range("J21").Value = 204
$ws.Range("K21").Value = -10.784313725490
$ws.Range("L21").Value = -18.75
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = -59.465478841870
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = 150
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 33.333333333333
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 44
$ws.Range("K23").Value = -22.727272727272
$ws.Range("L23").Value = 13.333333333333
$ws.Range("M23").Value = 3.030303030303
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -37.5
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -18.072289156626
$ws.Range("I24").Value = 228
$ws.Range("J24").Value = 243
$ws.Range("K24").Value = -6.172839506172
$ws.Range("L24").Value = 0.440528634361
$ws.Range("M24").Value = 72.727272727272
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = -29.310344827586
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 158
$ws.Range("K25").Value = -18.987341772151
$ws.Range("L25").Value = -3.759398496240
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -51.515151515151
$ws.Range("I26").Value = 79
$ws.Range("J26").Value = 98
$ws.Range("K26").Value = -19.387755102040
$ws.Range("L26").Value = -10.227272727272
$ws.Range("M26").Value = 31.666666666666
$ws.Range("G27").Value = 1
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = 120
$ws.Range("L28").Value = 57.142857142857
$ws.Range("N29").Value = -80
$ws.Range("N30").Value = -75
$ws.Range("G31").Value = 4
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = -80

# --- Numeric cells becoming text placeholders ("0" / "***.*") ---
# Step 1: copy value+type from a stable text anchor; Step 2: copy format from same anchor
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4104)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4104)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4104)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4104)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)

# --- Text placeholder cells becoming numeric values ---
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4104)
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D31").PasteSpecial(-4104)
$ws.Range("I14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 2
$ws.Range("L14").Copy()
$ws.Range("E31").PasteSpecial(-4104)
$ws.Range("L14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100

$excel.CutCopyMode = 0
